$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Shared-string content fix: JSON schema examples use type "int" -> "integer"
#    Affects every cell whose text is one of these two JSON snippets
#    (they occur on several "example" sheets sharing the same template text).
# ---------------------------------------------------------------------------
$oldSchema2 = '[{"name":"t0","type":"int"},
{"name":"t1","type":"long"}]'
$newSchema2 = '[{"name":"t0","type":"integer"},
{"name":"t1","type":"long"}]'

$oldSchema1 = '[{"name":"t0","type":"int"}]'
$newSchema1 = '[{"name":"t0","type":"integer"}]'

foreach ($sheetName in @("Template", "Folder", "File-timeseries", "File-relation")) {
    $s = $wb.Worksheets.Item($sheetName)
    foreach ($cellAddr in @("C2", "C3", "E2", "E3")) {
        $rng = $s.Range($cellAddr)
        if ($rng.Text -eq $oldSchema2) {
            $rng.Value = $newSchema2
        }
    }
}

$calcSheet = $wb.Worksheets.Item("File-calculate")
foreach ($cellAddr in @("D2", "D3")) {
    $rng = $calcSheet.Range($cellAddr)
    if ($rng.Text -eq $oldSchema1) {
        $rng.Value = $newSchema1
    }
}

# ---------------------------------------------------------------------------
# 2. Per-sheet selection / active-cell updates.
#    Order matters only for which sheet/cell ends up as the workbook's
#    overall active tab + selection: that must be File-reference / D24,
#    so it is activated last.
# ---------------------------------------------------------------------------
$selections = [ordered]@{
    "Template"          = "F14"
    "Folder"            = "E10"
    "File-timeseries"   = "D15"
    "File-relation"     = "D17"
    "File-calculate"    = "D19"
    "File-aggregation"  = "D20"
    "File-reference"    = "D24"
}

foreach ($name in $selections.Keys) {
    $sh = $wb.Worksheets.Item($name)
    [void]$sh.Activate()
    [void]$sh.Range($selections[$name]).Select()
}
